$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi re-run following Dr Hou's advice: recompute the Wnt1-Fzd3 LR-pair
# table with the new "ECs" sending/target cluster added alongside the
# existing "FAPs"/"sCs" clusters. Columns A:T are:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average/total expression value,
#   Ligand derived specificity of average/total expression value,
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average/total expression value,
#   Receptor derived specificity of average/total expression value,
#   Edge average/total expression weight,
#   Edge average/total expression derived specificity
$data = @(
    @("ECs", "Wnt1", "Fzd3", "ECs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 2, 0.6666666666666666, 0.129292, 0.387876, 0.06615700391713267, 0.06615700391713268, 0.003506786915999999, 0.031561082244, 0.004717950174570628, 0.004717950174570631),
    @("ECs", "Wnt1", "Fzd3", "FAPs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 0.4307096666666667, 1.292129, 0.220388431649395, 0.220388431649395, 0.011682138289, 0.105139244601, 0.01571687921170109, 0.0157168792117011),
    @("ECs", "Wnt1", "Fzd3", "sCs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1, 1.394319, 4.182957, 0.7134545644334723, 0.7134545644334724, 0.037818114237, 0.340363028133, 0.05087961799227442, 0.05087961799227444),
    @("FAPs", "Wnt1", "Fzd3", "ECs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 2, 0.6666666666666666, 0.129292, 0.387876, 0.06615700391713267, 0.06615700391713268, 0.04566679634666666, 0.41100116712, 0.06143905374256203, 0.06143905374256205),
    @("FAPs", "Wnt1", "Fzd3", "FAPs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 0.4307096666666667, 1.292129, 0.220388431649395, 0.220388431649395, 0.1521295256644445, 1.36916573098, 0.2046715524376939, 0.2046715524376939),
    @("FAPs", "Wnt1", "Fzd3", "sCs", 3, 1, 0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1, 1.394319, 4.182957, 0.7134545644334723, 0.7134545644334724, 0.49248276626, 4.43234489634, 0.6625749464411979, 0.662574946441198)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
